$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 2 ("Module 1, 2, 3: ICEES") - append text to the run that
# currently reads "(using N X N, feature_association2?)"
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(1)
$tr2 = $shp2.TextFrame.TextRange

$needle2 = "(using N X N, feature_association2?)"
$fullText2 = $tr2.Text
$idx2 = $fullText2.IndexOf($needle2)
if ($idx2 -ge 0) {
    $sub2 = $tr2.Characters($idx2 + 1, $needle2.Length)
    $sub2.Text = "(using N X N, feature_association2? Possibly using /features endpoint…)"
}

# ---------------------------------------------------------------
# Slide 3 ("Module 4: ...") - resize/reposition the rounded
# rectangle, update its text, and remove the four connector
# arrows + four rotated textbox labels that used to sit on top
# of it.
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(1)

# NOTE: Shape position/size is stored as single-precision (32-bit)
# points internally and truncated (not rounded) back to EMU on save,
# so a plain "EMU / 12700" division can land one EMU short of the
# target. The literals below are the float32 point values that
# truncate to the exact target EMU (1929384, 1602522, 7999325,
# 3285893) once re-expanded by PowerPoint's 12700 EMU/pt.
$shp3.Left = 151.92001342773438
$shp3.Top = 126.1828384399414
$shp3.Width = 629.8681640625
$shp3.Height = 258.73175048828125

$tr3 = $shp3.TextFrame.TextRange
$needle3 = "Module 4: New Tools/Analysis? (possibly Gamma?) "
$fullText3 = $tr3.Text
$idx3 = $fullText3.IndexOf($needle3)
if ($idx3 -ge 0) {
    $sub3 = $tr3.Characters($idx3 + 1, $needle3.Length)
    $sub3.Text = "Module 4: New Tools/Analysis? (possibly Gamma is capable of this already?) "
}

# Remove the trailing shapes (4 straight arrow connectors + 4
# rotated textboxes) in reverse order, leaving only the rounded
# rectangle (Shape 1) behind.
while ($s3.Shapes.Count -gt 1) {
    $s3.Shapes.Item($s3.Shapes.Count).Delete()
}
